$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Schors/ wankant"
$ws.Range("F1").Value = "Ontbrekend spint"
$ws.Range("D1").Value = "Spinthout- ringen"

$ws.Range("D5").Select()
